# Insert a new weekly price record at row 56 ("Hortaliza, Terminal La
# Palmera de La Serena - Berenjena"). This pushes the existing rows
# 56..82 down to 57..83 (dimension grows from A1:R82 to A1:R83) and
# populates the newly opened row 56 with the new data point.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 56:82 down to 57:83, opening up a blank row 56.
$ws.Rows("56:56").Insert()

# Populate the new row 56 with the inserted record's data.
$ws.Range("A56").Value = 8
$ws.Range("B56").Value = "Terminal La Palmera de La Serena"
$ws.Range("C56").Value = "Coquimbo"
$ws.Range("D56").Value = 44466
$ws.Range("E56").Value = 4
$ws.Range("F56").Value = 100112001
$ws.Range("G56").Value = "Berenjena"
$ws.Range("H56").Value = "Sin especificar"
$ws.Range("I56").Value = "Primera"
$ws.Range("J56").Value = 600
$ws.Range("K56").Value = 8500
$ws.Range("L56").Value = 9000
$ws.Range("M56").Value = 8750
$ws.Range("N56").Value = "$/caja 60 unidades"
$ws.Range("O56").Value = "Región de Arica y Parinacota"
$ws.Range("P56").Value = 146
$ws.Range("Q56").Value = 60
$ws.Range("R56").Value = "Hortaliza"
